# Insert 9 new daily-history rows (2019-11-18 .. 2019-11-28) right before the
# existing 2019-11-29 row, shifting all subsequent rows down by 9 (old row 605
# -> new row 614, ..., old row 677 -> new row 686). Dimension grows from
# A1:I677 to A1:I686.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: make room for the new rows -------------------------------------
# Old rows 605..613 (and everything after) need to shift down by 9 rows, so
# insert 9 blank rows starting at row 605.
$ws.Range("A605:A613").EntireRow.Insert()

# --- Step 2: populate the freshly inserted rows ------------------------------
$newRows = @(
    @{ Row = 605; Ts = 1574035200; Date = "2019-11-18"; Open = 0.275; High = 0.275; Low = 0.27;  Close = 0.27;  Vol = 852600 },
    @{ Row = 606; Ts = 1574121600; Date = "2019-11-19"; Open = 0.27;  High = 0.3;   Low = 0.27;  Close = 0.3;   Vol = 14486600 },
    @{ Row = 607; Ts = 1574208000; Date = "2019-11-20"; Open = 0.3;   High = 0.3;   Low = 0.285; Close = 0.285; Vol = 4057600 },
    @{ Row = 608; Ts = 1574294400; Date = "2019-11-21"; Open = 0.285; High = 0.29;  Low = 0.285; Close = 0.29;  Vol = 1375400 },
    @{ Row = 609; Ts = 1574380800; Date = "2019-11-22"; Open = 0.29;  High = 0.29;  Low = 0.28;  Close = 0.285; Vol = 3052400 },
    @{ Row = 610; Ts = 1574640000; Date = "2019-11-25"; Open = 0.29;  High = 0.29;  Low = 0.275; Close = 0.285; Vol = 1851700 },
    @{ Row = 611; Ts = 1574726400; Date = "2019-11-26"; Open = 0.28;  High = 0.28;  Low = 0.275; Close = 0.28;  Vol = 2607600 },
    @{ Row = 612; Ts = 1574812800; Date = "2019-11-27"; Open = 0.275; High = 0.28;  Low = 0.27;  Close = 0.28;  Vol = 1374800 },
    @{ Row = 613; Ts = 1574899200; Date = "2019-11-28"; Open = 0.275; High = 0.275; Low = 0.27;  Close = 0.275; Vol = 593800 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Range("A$row").Value = $r.Ts

    # Force the text-like columns (date / id / name) to stay text instead of
    # being auto-coerced into a date serial number or a plain number, then
    # clear the style back to Normal so no stray number-format style sticks.
    $ws.Range("B$row").NumberFormat = "@"
    $ws.Range("B$row").Value = $r.Date
    $ws.Range("B$row").Style = "Normal"

    $ws.Range("C$row").NumberFormat = "@"
    $ws.Range("C$row").Value = "0192"
    $ws.Range("C$row").Style = "Normal"

    $ws.Range("D$row").NumberFormat = "@"
    $ws.Range("D$row").Value = "INTA"
    $ws.Range("D$row").Style = "Normal"

    $ws.Range("E$row").Value = $r.Open
    $ws.Range("F$row").Value = $r.High
    $ws.Range("G$row").Value = $r.Low
    $ws.Range("H$row").Value = $r.Close
    $ws.Range("I$row").Value = $r.Vol
}
